# Reorder the product rows of the "Feuil1" sheet (rows 2-11, columns A-E).
# The underlying data does not change - only the row order - so every row
# is read into memory first, and then the rows are written back out in the
# new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$firstRow = 2
$lastRow = 11
$lastCol = 5

# Capture the current rows (as an array of row-arrays) before overwriting anything.
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $row += ,$ws.Cells.Item($r, $c).Value2
    }
    $rows += ,$row
}

# New row order, expressed as the 0-based index (within $rows) that should
# land on each destination row (2..11, in order). Product codes 555, 222,
# 777, 111, 919, 888, 444, 999, 666, 333 (top to bottom) after the edit.
$order = @(9, 5, 0, 4, 7, 8, 2, 6, 1, 3)

for ($i = 0; $i -lt $order.Count; $i++) {
    $srcRow = $rows[$order[$i]]
    $destRow = $firstRow + $i
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcRow[$c - 1]
    }
}
